$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force text to avoid numeric auto-conversion, then restore default (no) style ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.694.59'
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.323.38'
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.94'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.78'
$ws.Range("D6").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.319.89'
$ws.Range("D9").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.41'
$ws.Range("D12").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '667.02'
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.859.66'
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.37'
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.859.51'
$ws.Range("D17").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.319.52'
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.42'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.89'
$ws.Range("D21").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.98'
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.47'
$ws.Range("D25").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.27'
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.39'
$ws.Range("D29").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.31'
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '585.55'
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.93'
$ws.Range("D33").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.714.28'
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.75'
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.25'
$ws.Range("D38").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '32.69'
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.62'
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.10'
$ws.Range("D42").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0661'
$ws.Range("D44").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0406'
$ws.Range("D46").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.128'
$ws.Range("D48").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '126.54'
$ws.Range("D51").ClearFormats()

# --- Volume(1h) column (E): plain text assignment (percent strings are never numeric-parseable) ---
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  -2.19%  '
$ws.Range("E6").Value = '  -7.12%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -2.78%  '
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("E10").Value = '  -4.61%  '
$ws.Range("E11").Value = '  -2.54%  '
$ws.Range("E12").Value = '  -4.44%  '
$ws.Range("E13").Value = '  -4.03%  '
$ws.Range("E14").Value = '  +4.44%  '
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("E16").Value = '  -2.93%  '
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("E19").Value = '  -1.69%  '
$ws.Range("E20").Value = '  -3.74%  '
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("E22").Value = '  -2.66%  '
$ws.Range("E23").Value = '  +5.40%  '
$ws.Range("E24").Value = '  -5.33%  '
$ws.Range("E25").Value = '  -2.33%  '
$ws.Range("E26").Value = '  -5.12%  '
$ws.Range("E27").Value = '  -6.95%  '
$ws.Range("E28").Value = '  -4.84%  '
$ws.Range("E29").Value = '  +1.23%  '
$ws.Range("E31").Value = '  +5.79%  '
$ws.Range("E32").Value = '  -4.44%  '
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("E34").Value = '  -2.44%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").Value = '  -8.42%  '
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("E38").Value = '  -14.56%  '
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("E40").Value = '  -3.36%  '
$ws.Range("E41").Value = '  -6.06%  '
$ws.Range("E42").Value = '  -4.94%  '
$ws.Range("E43").Value = '  -3.53%  '
$ws.Range("E44").Value = '  -5.65%  '
$ws.Range("E45").Value = '  -5.17%  '
$ws.Range("E46").Value = '  -3.89%  '
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("E48").Value = '  -2.04%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("E50").Value = '  -4.03%  '
$ws.Range("E51").Value = '  -1.20%  '
